$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.563.74"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +5.07%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.843.00"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  +4.19%  "

$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.028"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  +2.35%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.99"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +4.61%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.024"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +2.20%  "

$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4378"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +3.80%  "

$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3739"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +4.10%  "

$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07390"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +3.88%  "

$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8767"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +4.99%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.53"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +6.19%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.859.49"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +5.35%  "

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.495"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +5.29%  "

$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.682"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +3.92%  "

$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07146"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +3.91%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.73"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +5.28%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.029"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +2.61%  "

$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009005"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +3.96%  "

$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.024"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +2.24%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.42"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +3.51%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.566.88"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +4.82%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.253"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +3.26%  "

$ws.Range("E23").Value = "  +2.38%  "

$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.075.14"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +4.76%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.96"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +3.40%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.921"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +7.26%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.73"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +4.31%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.271"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +4.55%  "

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.940"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +5.90%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.30"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +2.15%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09089"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +2.99%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.209"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +8.99%  "

$ws.Range("E33").Value = "  +6.23%  "

$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.504"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +5.00%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.866"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +4.83%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.026"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +2.55%  "

$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.144"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  +3.11%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01974"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +5.10%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05257"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +3.25%  "

$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5177"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +5.81%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.794"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +7.90%  "

$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1669"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  +4.11%  "

$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.653"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +5.31%  "

$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.537"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +6.83%  "

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "108.98"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +4.46%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.61"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +4.13%  "

$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.026"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +2.48%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.708"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +6.11%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4652"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +4.72%  "

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.909"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +13.48%  "

$ws.Range("E51").Value = "  +2.99%  "

